# Apply the "comma-separated -> semicolon-separated" list wording/value edit
# described in the commit message:
#   "comma-separated lists in input replaced with semicolon-separated lists"

$wb = $excel.ActiveWorkbook

# Unicode curly single quotes used in the original description text.
$lsq = [char]0x2018   # left single quotation mark
$rsq = [char]0x2019   # right single quotation mark

# --- Sheet "Description" ---
$wsDesc = $wb.Worksheets.Item("Description")

# Row 4 (ChromosomeEnds) - Values column description
$wsDesc.Range("C4").Value = "Semicolon separated list of integers, or # to leave empty. Must be # for haploid (asexual) systems."

# Row 11 (PatchList) - Values column description
$wsDesc.Range("C11").Value = "Either 1) Semicolon-separated list (CSL), 2) " + $lsq + "random" + $rsq + ", 3) " + $lsq + "random_occupied" + $rsq + " or 4) " + $lsq + "all" + $rsq + ". Patch 0 is reserved for an internal matrix and cannot be selected (if the CSL option is chosen). In the CSL and random case, there is no internal check of whether the specified/sampled patch exist through the simulation, so it is up to the user to ensure this is the case, and special care should be given for dynamic landscapes."

# Row 14 (Stages) - Values column description
$wsDesc.Range("C14").Value = 'Semicolon-separated list or "all"'

# --- Sheet "GeneticsFile" ---
$wsGen = $wb.Worksheets.Item("GeneticsFile")

# Row 3 sample data: comma separated lists -> semicolon separated lists
$wsGen.Range("C3").Value = "5;10;15"
$wsGen.Range("J3").Value = "4;5;8"
$wsGen.Range("M3").Value = "1;2"

# Restore cursor/selection state: leave "GeneticsFile" with M3 selected,
# then activate "Description" (the saved-as-active tab) with C12 selected.
$wsGen.Activate() | Out-Null
$wsGen.Range("M3").Select() | Out-Null

$wsDesc.Activate() | Out-Null
$wsDesc.Range("C12").Select() | Out-Null
